$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 33 with the same pattern as the preceding data rows
$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 10032
$ws.Range("C33").Value = "eng"
$ws.Range("D33").Value = $true
$ws.Range("E33").Value = "superadmin"
$ws.Range("F33").Value = "now()"
$ws.Range("G33").Value = "now()"

# Update the selection to match the state captured after the edit
$ws.Range("B30").Select()
